$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- CasesTab row (row 2): append an ORDER BY / LIMIT clause to the Cypher query ---
$b2 = $ws.Range("B2").Value2
$ws.Range("B2").Value2 = $b2 + "`n order By ss.study_subject_id ASC LIMIT 100 "

# --- SamplesTab row (row 3): append an ORDER BY / LIMIT clause to the Cypher query ---
$b3 = $ws.Range("B3").Value2
$ws.Range("B3").Value2 = $b3 + "`n order By samp.sample_id ASC LIMIT 100"

# --- FilesTab row (row 4): replace the trailing "order by f.file_name" with the new clause ---
$b4 = $ws.Range("B4").Value2
$oldSuffix = "order by f.file_name"
$b4 = $b4.Substring(0, $b4.Length - $oldSuffix.Length) + "order By f.file_name ASC LIMIT 100"
$ws.Range("B4").Value2 = $b4

# Row heights grew because the wrapped query text got longer
$ws.Rows.Item(2).RowHeight = 331.2
$ws.Rows.Item(3).RowHeight = 360

# Update selection / scroll position recorded in the sheet view
$ws.Range("C3").Select() | Out-Null
